$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target final state for rows 2-7 (title | timestamp | historical distance | time bucket | uri)
$ws.Range("A2").Value = "Coast Chinatown Seeks to Still Fears"
$ws.Range("B2").Value = "1977-10-02T01:00:00UTC"
$ws.Range("C2").Value = 28
$ws.Range("D2").Value = "day_2_to_30"
$ws.Range("E2").Value = "https://www.nytimes.com/1977/10/02/archives/coast-chinatown-seeks-to-still-fears.html"

$ws.Range("A3").Value = '$25,000 Reward Offered In 5 Slayings on Coast'
$ws.Range("B3").Value = "1977-09-07T01:00:00UTC"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "day_2_to_30"
$ws.Range("E3").Value = "https://www.nytimes.com/1977/09/07/archives/25000-reward-offered-in-5-slayings-on-coast.html"

$ws.Range("A4").Value = "Around the Nation"
$ws.Range("B4").Value = "1978-04-22T00:00:00UTC"
$ws.Range("C4").Value = 230
$ws.Range("D4").Value = "day_31_beyond"
$ws.Range("E4").Value = "https://www.nytimes.com/1978/04/22/archives/around-the-nation-two-acquitted-in-georgia-despite-carter-testimony.html"

$ws.Range("A5").Value = "Year of Horse a Time of Change For Chinatowns Across Nation"
$ws.Range("B5").Value = "1978-02-13T00:00:00UTC"
$ws.Range("C5").Value = 162
$ws.Range("D5").Value = "day_31_beyond"
$ws.Range("E5").Value = "https://www.nytimes.com/1978/02/13/archives/year-of-horse-a-time-of-change-for-chinatowns-across-nation-i-want.html"

$ws.Range("A6").Value = "SQUAD REDUCES CHINATOWN VIOLENCE ON COAST"
$ws.Range("B6").Value = "1-01-01T00:00:00UTC"
$ws.Range("C6").Value = "unknown"
$ws.Range("D6").Value = "unknown"
$ws.Range("E6").Value = "https://www.nytimes.com/1983/09/23/us/squad-reduces-chinatown-violence-on-coast.html"

$ws.Range("A7").Value = "The Golden Dragon Restaurant Massacre"
$ws.Range("B7").Value = "1-01-01T00:00:00UTC"
$ws.Range("C7").Value = "unknown"
$ws.Range("D7").Value = "unknown"
$ws.Range("E7").Value = "http://foundsf.org/index.php?title=The_Golden_Dragon_Restaurant_Massacre"

# Rebuild the hyperlinks so their targets line up with the new row order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.nytimes.com/1977/10/02/archives/coast-chinatown-seeks-to-still-fears.html")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.nytimes.com/1977/09/07/archives/25000-reward-offered-in-5-slayings-on-coast.html")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.nytimes.com/1978/04/22/archives/around-the-nation-two-acquitted-in-georgia-despite-carter-testimony.html")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.nytimes.com/1978/02/13/archives/year-of-horse-a-time-of-change-for-chinatowns-across-nation-i-want.html")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.nytimes.com/1983/09/23/us/squad-reduces-chinatown-violence-on-coast.html")
$ws.Hyperlinks.Add($ws.Range("E7"), "http://foundsf.org/index.php?title=The_Golden_Dragon_Restaurant_Massacre")
